$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- Sheet ALC ---
# Row 33
$ws_ALC.Range("H33").Value = 320.44
$ws_ALC.Range("I33").Value = 344.13635
$ws_ALC.Range("K33").Value = 344.13635
$ws_ALC.Range("M33").Value = -115.13635

# Row 130
$ws_ALC.Range("H130").Value = 58182.5
$ws_ALC.Range("J130").Value = 58182.5
$ws_ALC.Range("L130").Value = 58182.5
$ws_ALC.Range("N130").Value = -68222.5

# Row 132
$ws_ALC.Range("H132").Value = 2161.1719
$ws_ALC.Range("I132").Value = 1966.356
$ws_ALC.Range("K132").Value = 5899.068
$ws_ALC.Range("M132").Value = -3369.068

# Row 137
$ws_ALC.Range("H137").Value = 5128959
$ws_ALC.Range("I137").Value = 674.375
$ws_ALC.Range("J137").Value = 13334215
$ws_ALC.Range("K137").Value = 2023.125
$ws_ALC.Range("L137").Value = 40002645
$ws_ALC.Range("M137").Value = 526.875
$ws_ALC.Range("N137").Value = -40007745

# --- Sheet ARM ---
# Row 32
$ws_ARM.Range("H32").Value = 8169.348
$ws_ARM.Range("I32").Value = 8449.457
$ws_ARM.Range("J32").Value = 7137.3687
$ws_ARM.Range("K32").Value = 8449.457
$ws_ARM.Range("L32").Value = 7137.3687
$ws_ARM.Range("M32").Value = -8162.457
$ws_ARM.Range("N32").Value = -7711.3687

# Row 61
$ws_ARM.Range("H61").Value = 11906322
$ws_ARM.Range("I61").Value = 13515159
$ws_ARM.Range("J61").Value = 928.4
$ws_ARM.Range("K61").Value = 13515159
$ws_ARM.Range("L61").Value = 928.4
$ws_ARM.Range("M61").Value = -13514947
$ws_ARM.Range("N61").Value = -1352.4

# Row 74
$ws_ARM.Range("H74").Value = 10418333
$ws_ARM.Range("I74").Value = 13159475
$ws_ARM.Range("K74").Value = 13159475
$ws_ARM.Range("M74").Value = -13158601

# Row 77
$ws_ARM.Range("H77").Value = 10418333
$ws_ARM.Range("I77").Value = 13159475
$ws_ARM.Range("K77").Value = 65797375
$ws_ARM.Range("M77").Value = -65793007

# Row 132
$ws_ARM.Range("H132").Value = 5210105.5
$ws_ARM.Range("I132").Value = 6251503.5
$ws_ARM.Range("J132").Value = 3117.75
$ws_ARM.Range("K132").Value = 18754510.5
$ws_ARM.Range("L132").Value = 9353.25
$ws_ARM.Range("M132").Value = -18751980.5
$ws_ARM.Range("N132").Value = -14413.25

# Row 136
$ws_ARM.Range("H136").Value = 11906322
$ws_ARM.Range("I136").Value = 13515159
$ws_ARM.Range("J136").Value = 928.4
$ws_ARM.Range("K136").Value = 40545477
$ws_ARM.Range("L136").Value = 2785.2
$ws_ARM.Range("M136").Value = -40542927
$ws_ARM.Range("N136").Value = -7885.2

# --- Sheet BSM ---
# Row 7
$ws_BSM.Range("H7").Value = 22060.6
$ws_BSM.Range("I7").Value = 33934.332
$ws_BSM.Range("J7").Value = 4250
$ws_BSM.Range("K7").Value = 33934.332
$ws_BSM.Range("L7").Value = 4250
$ws_BSM.Range("M7").Value = -33821.332
$ws_BSM.Range("N7").Value = -4476

# Row 134
$ws_BSM.Range("H134").Value = 2429.4644
$ws_BSM.Range("I134").Value = 1550.1904
$ws_BSM.Range("J134").Value = 5067.2856
$ws_BSM.Range("K134").Value = 4650.5712
$ws_BSM.Range("L134").Value = 15201.8568
$ws_BSM.Range("M134").Value = -2115.5712
$ws_BSM.Range("N134").Value = -20271.8568

# --- Sheet CRP ---
# Row 31
$ws_CRP.Range("H31").Value = 4834315.5
$ws_CRP.Range("I31").Value = 3294.0688
$ws_CRP.Range("K31").Value = 3294.0688
$ws_CRP.Range("M31").Value = -2999.0688

# Row 34
$ws_CRP.Range("H34").Value = 4834315.5
$ws_CRP.Range("I34").Value = 3294.0688
$ws_CRP.Range("K34").Value = 3294.0688
$ws_CRP.Range("M34").Value = -3092.0688

# Row 58
$ws_CRP.Range("H58").Value = 1898.5555
$ws_CRP.Range("I58").Value = 1081.125
$ws_CRP.Range("K58").Value = 1081.125
$ws_CRP.Range("M58").Value = -878.125

# Row 60
$ws_CRP.Range("H60").Value = 8029.9
$ws_CRP.Range("J60").Value = 8549.833000000001
$ws_CRP.Range("L60").Value = 8549.833000000001
$ws_CRP.Range("N60").Value = -9571.833000000001

# Row 132
$ws_CRP.Range("H132").Value = 8065796.5
$ws_CRP.Range("I132").Value = 9260273
$ws_CRP.Range("J132").Value = 3076.5
$ws_CRP.Range("K132").Value = 27780819
$ws_CRP.Range("L132").Value = 9229.5
$ws_CRP.Range("M132").Value = -27778289
$ws_CRP.Range("N132").Value = -14289.5

# Row 134
$ws_CRP.Range("H134").Value = 1524.1163
$ws_CRP.Range("I134").Value = 1395.925
$ws_CRP.Range("J134").Value = 3233.3333
$ws_CRP.Range("K134").Value = 4187.775
$ws_CRP.Range("L134").Value = 9699.999899999999
$ws_CRP.Range("M134").Value = -1652.775
$ws_CRP.Range("N134").Value = -14769.9999

# Row 136
$ws_CRP.Range("H136").Value = 1898.5555
$ws_CRP.Range("I136").Value = 1081.125
$ws_CRP.Range("K136").Value = 3243.375
$ws_CRP.Range("M136").Value = -693.375

# --- Sheet CUL ---
# Row 34
$ws_CUL.Range("H34").Value = 1317.7142
$ws_CUL.Range("I34").Value = 181.25
$ws_CUL.Range("K34").Value = 543.75
$ws_CUL.Range("M34").Value = -459.75

# Row 130
$ws_CUL.Range("H130").Value = 5435
$ws_CUL.Range("J130").Value = 6146.154
$ws_CUL.Range("L130").Value = 18438.462
$ws_CUL.Range("N130").Value = -28478.462

# Row 136
$ws_CUL.Range("H136").Value = 3039
$ws_CUL.Range("I136").Value = 1826
$ws_CUL.Range("J136").Value = 3544.4167
$ws_CUL.Range("K136").Value = 5478
$ws_CUL.Range("L136").Value = 10633.2501
$ws_CUL.Range("M136").Value = -378
$ws_CUL.Range("N136").Value = -20833.2501

# Row 139
$ws_CUL.Range("H139").Value = 2967.4285
$ws_CUL.Range("I139").Value = 1343.3334
$ws_CUL.Range("J139").Value = 7027.6665
$ws_CUL.Range("K139").Value = 4030.0002
$ws_CUL.Range("L139").Value = 21082.9995
$ws_CUL.Range("M139").Value = 1109.9998
$ws_CUL.Range("N139").Value = -31362.9995

# Row 140
$ws_CUL.Range("H140").Value = 3418.4211
$ws_CUL.Range("I140").Value = 1639.2858
$ws_CUL.Range("J140").Value = 8400
$ws_CUL.Range("K140").Value = 4917.857400000001
$ws_CUL.Range("L140").Value = 25200
$ws_CUL.Range("M140").Value = 262.1425999999992
$ws_CUL.Range("N140").Value = -35560

# --- Sheet GSM ---
# Row 58
$ws_GSM.Range("H58").Value = 10046
$ws_GSM.Range("J58").Value = 10046
$ws_GSM.Range("L58").Value = 10046
$ws_GSM.Range("N58").Value = -10600

# --- Sheet LTW ---
# Row 57
$ws_LTW.Range("H57").Value = 20000
$ws_LTW.Range("J57").Value = 20000
$ws_LTW.Range("L57").Value = 20000
$ws_LTW.Range("N57").Value = -21132

# Row 58
$ws_LTW.Range("H58").Value = 0
$ws_LTW.Range("I58").Value = 0
$ws_LTW.Range("K58").Value = 0

# Row 75
$ws_LTW.Range("H75").Value = 0
$ws_LTW.Range("I75").Value = 0
$ws_LTW.Range("K75").Value = 0

# Row 78
$ws_LTW.Range("H78").Value = 0
$ws_LTW.Range("I78").Value = 0
$ws_LTW.Range("K78").Value = 0

# Row 132
$ws_LTW.Range("H132").Value = 6855204.5
$ws_LTW.Range("I132").Value = 3446.1091
$ws_LTW.Range("J132").Value = 27791134
$ws_LTW.Range("K132").Value = 10338.3273
$ws_LTW.Range("L132").Value = 83373402
$ws_LTW.Range("M132").Value = -7808.327300000001
$ws_LTW.Range("N132").Value = -83378462

# Row 136
$ws_LTW.Range("H136").Value = 10872985
$ws_LTW.Range("I136").Value = 12196232
$ws_LTW.Range("J136").Value = 22361
$ws_LTW.Range("K136").Value = 36588696
$ws_LTW.Range("L136").Value = 67083
$ws_LTW.Range("M136").Value = -36586146
$ws_LTW.Range("N136").Value = -72183

# Row 140
$ws_LTW.Range("H140").Value = 49641.875
$ws_LTW.Range("J140").Value = 49641.875
$ws_LTW.Range("L140").Value = 49641.875
$ws_LTW.Range("N140").Value = -60001.875

# --- Sheet WVR ---
# Row 136
$ws_WVR.Range("H136").Value = 843.6222
$ws_WVR.Range("I136").Value = 636.7805
$ws_WVR.Range("J136").Value = 2963.75
$ws_WVR.Range("K136").Value = 1910.3415
$ws_WVR.Range("L136").Value = 8891.25
$ws_WVR.Range("M136").Value = 639.6585
$ws_WVR.Range("N136").Value = -13991.25

# --- Cell removals (values present in before, absent in after) ---
$ws_LTW.Range("M58").ClearContents()
$ws_LTW.Range("M75").ClearContents()
$ws_LTW.Range("M78").ClearContents()
